# Apply the authored edit to the "jobs" worksheet:
#  - Row 2 is replaced with a new job listing (KAM - D2C Freight Operations / Weekday)
#  - Rows 3-6 (the other four job listings) are removed entirely
#  - The used range / dimension shrinks from A1:K6 to A1:K2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Overwrite row 2 with the new job listing's data ---
$ws.Range("A2").Value = "KAM - D2C Freight Operations"
$ws.Range("B2").Value = "Chandigarh, Chandigarh, India"
$ws.Range("C2").Value = "Weekday"
$ws.Range("D2").Value = "A providing logistics and D2C freight management services."
$ws.Range("E2").Value = "https://kg.diffbot.com/image/api/get?fetch=yes&urls=g%3Cj7P0Stn8p.%5DjEp9G.Bd%7Bk%3ESQLnc%7E%5B-AYYy-Z%3C%5BL-D%3A%7Bj-a%5CErs%3AwPm7S%5CzN_%3Cb4E.A%5Dm,g%3Cj7P0Stn8p.OgAu%3CRsCX3_1BY%2F%5Dt.%5B8LwifzFeCd.Nr0QoN%5CtW4BmCRh7%5E%7CZ-Oh1Cy6t%3Dm5Og.Z%5Et%60,g%3Cj7P0SttAc%3ANu.d9U%7CB%5DpM08.%7DZnIa6%40_%3Bv2FpDLhawxEr-%5BuTpJZ%40-r%3Ec2Jh9e%3ERjWglef5jxk%5DBa4MY%2Ftg"
$ws.Range("F2").Value = "5+ YOE2+ MgmtEscalation management, client liaison, dashboards, NDR, COD remittance; 5+ years in D2C freight ops; leadership/account management experience; strong logistics knowledge."
$ws.Range("G2").Value = "Microsoft Excel"
$ws.Range("H2").Value = "5+ YOE"
$ws.Range("I2").Value = "Full Time"
$ws.Range("J2").Value = "Onsite"
$ws.Range("K2").Value = "https://apply.workable.com/weekday-1/j/C7AA0742E4/"

# --- Remove the now-obsolete rows 3 through 6 entirely, shrinking the sheet to A1:K2 ---
$ws.Range("A3:K6").EntireRow.Delete() | Out-Null
